$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 42 (ALC)
$ws.Range("H42").Value = 483.5
$ws.Range("I42").Value = 448.75
$ws.Range("J42").Value = 553
$ws.Range("K42").Value = 1346.25
$ws.Range("L42").Value = 1659
$ws.Range("M42").Value = -1116.25
$ws.Range("N42").Value = -2119

# Row 86 (ALC)
$ws.Range("H86").Value = 2150.75
$ws.Range("I86").Value = 2150.75
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2150.75
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -1027.75

# Row 89 (ALC)
$ws.Range("H89").Value = 2150.75
$ws.Range("I89").Value = 2150.75
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10753.75
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -5137.75

# Row 137 (ALC)
$ws.Range("H137").Value = 1673.9
$ws.Range("I137").Value = 1205.4762
$ws.Range("J137").Value = 2766.889
$ws.Range("K137").Value = 3616.4286
$ws.Range("L137").Value = 8300.667000000001
$ws.Range("M137").Value = -1066.4286
$ws.Range("N137").Value = -13400.667

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 4793.608
$ws.Range("I32").Value = 3644.0667
$ws.Range("J32").Value = 13415.167
$ws.Range("K32").Value = 3644.0667
$ws.Range("L32").Value = 13415.167
$ws.Range("M32").Value = -3357.0667
$ws.Range("N32").Value = -13989.167

# Row 45 (ARM)
$ws.Range("H45").Value = 2191.5
$ws.Range("I45").Value = 2218.8572
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 2218.8572
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1841.8572
$ws.Range("N45").Value = -2754

# Row 122 (ARM)
$ws.Range("H122").Value = 2911.5
$ws.Range("I122").Value = 2905.6
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8716.799999999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6266.799999999999
$ws.Range("N122").Value = -13900

# Row 132 (ARM)
$ws.Range("H132").Value = 2767.4075
$ws.Range("I132").Value = 2073.5557
$ws.Range("J132").Value = 4155.1113
$ws.Range("K132").Value = 6220.6671
$ws.Range("L132").Value = 12465.3339
$ws.Range("M132").Value = -3690.6671
$ws.Range("N132").Value = -17525.3339

# Row 134 (ARM)
$ws.Range("H134").Value = 89000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 89000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 89000
$ws.Range("N134").Value = -99140

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (BSM)
$ws.Range("H86").Value = 7750
$ws.Range("I86").Value = 5500
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 5500
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -4377
$ws.Range("N86").Value = -12246

# Row 89 (BSM)
$ws.Range("H89").Value = 7750
$ws.Range("I89").Value = 5500
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 27500
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -21884
$ws.Range("N89").Value = -61232

# Row 119 (BSM)
$ws.Range("H119").Value = 60700
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 60700
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 60700
$ws.Range("N119").Value = -70376

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 949.6667
$ws.Range("I31").Value = 925
$ws.Range("J31").Value = 999
$ws.Range("K31").Value = 925
$ws.Range("L31").Value = 999
$ws.Range("M31").Value = -630
$ws.Range("N31").Value = -1589

# Row 34 (CRP)
$ws.Range("H34").Value = 949.6667
$ws.Range("I34").Value = 925
$ws.Range("J34").Value = 999
$ws.Range("K34").Value = 925
$ws.Range("L34").Value = 999
$ws.Range("M34").Value = -723
$ws.Range("N34").Value = -1403

# Row 58 (CRP)
$ws.Range("H58").Value = 2756.1538
$ws.Range("I58").Value = 2568.6
$ws.Range("J58").Value = 2873.375
$ws.Range("K58").Value = 2568.6
$ws.Range("L58").Value = 2873.375
$ws.Range("M58").Value = -2365.6
$ws.Range("N58").Value = -3279.375

# Row 107 (CRP)
$ws.Range("H107").Value = 2032.9
$ws.Range("I107").Value = 1121.75
$ws.Range("J107").Value = 2640.3333
$ws.Range("K107").Value = 1121.75
$ws.Range("L107").Value = 2640.3333
$ws.Range("M107").Value = 798.25
$ws.Range("N107").Value = -6480.3333

# Row 122 (CRP)
$ws.Range("H122").Value = 2950.1428
$ws.Range("I122").Value = 2209.5
$ws.Range("J122").Value = 3937.6667
$ws.Range("K122").Value = 6628.5
$ws.Range("L122").Value = 11813.0001
$ws.Range("M122").Value = -4178.5
$ws.Range("N122").Value = -16713.0001

# Row 132 (CRP)
$ws.Range("H132").Value = 2836.7144
$ws.Range("I132").Value = 2167.5715
$ws.Range("J132").Value = 4175
$ws.Range("K132").Value = 6502.7145
$ws.Range("L132").Value = 12525
$ws.Range("M132").Value = -3972.7145
$ws.Range("N132").Value = -17585

# Row 134 (CRP)
$ws.Range("H134").Value = 3651.6924
$ws.Range("I134").Value = 3547.9
$ws.Range("J134").Value = 3997.6667
$ws.Range("K134").Value = 10643.7
$ws.Range("L134").Value = 11993.0001
$ws.Range("M134").Value = -8108.700000000001
$ws.Range("N134").Value = -17063.0001

# Row 136 (CRP)
$ws.Range("H136").Value = 2756.1538
$ws.Range("I136").Value = 2568.6
$ws.Range("J136").Value = 2873.375
$ws.Range("K136").Value = 7705.799999999999
$ws.Range("L136").Value = 8620.125
$ws.Range("M136").Value = -5155.799999999999
$ws.Range("N136").Value = -13720.125

$ws = $wb.Worksheets.Item("CUL")
# Row 33 (CUL)
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()

# Row 64 (CUL)
$ws.Range("H64").Value = 750
$ws.Range("I64").Value = 750
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 2250
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0
$ws.Range("M64").Value = -1980

# Row 67 (CUL)
$ws.Range("H67").Value = 750
$ws.Range("I67").Value = 750
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 2250
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0
$ws.Range("M67").Value = -1314

# Row 113 (CUL)
$ws.Range("H113").Value = 1900.0526
$ws.Range("I113").Value = 2375
$ws.Range("J113").Value = 1844.1765
$ws.Range("K113").Value = 7125
$ws.Range("L113").Value = 5532.529500000001
$ws.Range("M113").Value = -4955
$ws.Range("N113").Value = -9872.529500000001

# Row 114 (CUL)
$ws.Range("H114").Value = 5949.5
$ws.Range("I114").Value = 5999
$ws.Range("J114").Value = 5900
$ws.Range("K114").Value = 17997
$ws.Range("L114").Value = 17700
$ws.Range("M114").Value = -14743
$ws.Range("N114").Value = -24208

# Row 137 (CUL)
$ws.Range("H137").Value = 3733.1667
$ws.Range("I137").Value = 3100
$ws.Range("J137").Value = 4999.5
$ws.Range("K137").Value = 9300
$ws.Range("L137").Value = 14998.5
$ws.Range("M137").Value = -4200
$ws.Range("N137").Value = -25198.5

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (GSM)
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = 0

# Row 102 (GSM)
$ws.Range("H102").Value = 1283.3572
$ws.Range("I102").Value = 1279.75
$ws.Range("J102").Value = 1305
$ws.Range("K102").Value = 1279.75
$ws.Range("L102").Value = 1305
$ws.Range("M102").Value = 342.25
$ws.Range("N102").Value = -4549

# Row 113 (GSM)
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

# Row 122 (GSM)
$ws.Range("H122").Value = 3249.75
$ws.Range("I122").Value = 2333
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 6999
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -4549
$ws.Range("N122").Value = -22900

# Row 126 (GSM)
$ws.Range("H126").Value = 1229.8
$ws.Range("I126").Value = 1037.25
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3111.75
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -641.75
$ws.Range("N126").Value = -10940

# Row 139 (GSM)
$ws.Range("H139").Value = 70713.8
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 70713.8
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 70713.8
$ws.Range("N139").Value = -80993.8

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 4000
$ws.Range("I7").Value = 4000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3888

# Row 55 (LTW)
$ws.Range("H55").Value = 225
$ws.Range("I55").Value = 250
$ws.Range("J55").Value = 200
$ws.Range("K55").Value = 250
$ws.Range("L55").Value = 200
$ws.Range("M55").Value = -77
$ws.Range("N55").Value = -546

# Row 82 (LTW)
$ws.Range("H82").Value = 1999.5
$ws.Range("I82").Value = 1999
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 1999
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -1638
$ws.Range("N82").Value = -2722

# Row 85 (LTW)
$ws.Range("H85").Value = 1999.5
$ws.Range("I85").Value = 1999
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 1999
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -751
$ws.Range("N85").Value = -4496

# Row 122 (LTW)
$ws.Range("H122").Value = 8000
$ws.Range("I122").Value = 8000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 24000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -21550

# Row 126 (LTW)
$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9530

# Row 136 (LTW)
$ws.Range("H136").Value = 19999750
$ws.Range("I136").Value = 19999750
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 59999250
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -59996700

$ws = $wb.Worksheets.Item("WVR")
# Row 113 (WVR)
$ws.Range("H113").Value = 1594.4
$ws.Range("I113").Value = 1531.375
$ws.Range("J113").Value = 1846.5
$ws.Range("K113").Value = 4594.125
$ws.Range("L113").Value = 5539.5
$ws.Range("M113").Value = -2424.125
$ws.Range("N113").Value = -9879.5

# Row 126 (WVR)
$ws.Range("H126").Value = 1999.2
$ws.Range("I126").Value = 1999.2
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5997.6
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3527.6

# Row 136 (WVR)
$ws.Range("H136").Value = 2016.6296
$ws.Range("I136").Value = 1691.3334
$ws.Range("J136").Value = 3155.1667
$ws.Range("K136").Value = 5074.0002
$ws.Range("L136").Value = 9465.500100000001
$ws.Range("M136").Value = -2524.0002
$ws.Range("N136").Value = -14565.5001

# Row 138 (WVR)
$ws.Range("H138").Value = 96300
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 96300
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 96300
$ws.Range("N138").Value = -106580
